# Update the "Census region" labels in column A so each level is prefixed
# with its numeric code, and refresh the dependent "Indicator" interaction
# labels to match (also prefixing the Indicator levels with 0./1.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value  = "Census region=1. NE"
$ws.Range("A8").Value  = "Census region=2. N Cntrl"
$ws.Range("A10").Value = "Census region=3. South"
$ws.Range("A12").Value = "Census region=4. West"

$ws.Range("A14").Value = "Census region=1. NE * Indicator=0. No"
$ws.Range("A16").Value = "Census region=1. NE * Indicator=1. Yes"
$ws.Range("A18").Value = "Census region=2. N Cntrl * Indicator=0. No"
$ws.Range("A20").Value = "Census region=2. N Cntrl * Indicator=1. Yes"
$ws.Range("A22").Value = "Census region=3. South * Indicator=0. No"
$ws.Range("A24").Value = "Census region=3. South * Indicator=1. Yes"
$ws.Range("A26").Value = "Census region=4. West * Indicator=0. No"
$ws.Range("A28").Value = "Census region=4. West * Indicator=1. Yes"
